# Update countries & provincias Spain
# - Swap Moldavia / Nigeria ordering (Moldavia's stats overtook Nigeria's)
# - Refresh the "Datos actualizados..." timestamp
# - Refresh case/death/recovered counters for several countries

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 18:05"

# --- Moldavia overtakes Nigeria: swap the two country rows ------------
# Row 61 was Nigeria / Row 62 was Moldavia -> now Row 61 = Moldavia (updated
# figures), Row 62 = Nigeria (keeps its previous figures).
$ws.Range("A61").Value = "Moldavia"
$ws.Range("A62").Value = "Nigeria"

# --- Numeric refreshes --------------------------------------------------
# Row 4 (Estados Unidos)
$ws.Range("B4").Value = 1599999
$ws.Range("C4").Value = 7276
$ws.Range("D4").Value = 371279
$ws.Range("E4").Value = 1133495
$ws.Range("G4").Value = 289
$ws.Range("H4").Value = 95225

# Row 11
$ws.Range("B11").Value = 178797
$ws.Range("C11").Value = 266
$ws.Range("E11").Value = 12526

# Row 35
$ws.Range("B35").Value = 20143
$ws.Range("C35").Value = 404
$ws.Range("E35").Value = 10719
$ws.Range("G35").Value = 10
$ws.Range("H35").Value = 972

# Row 45
$ws.Range("D45").Value = 7366
$ws.Range("E45").Value = 5843

# Row 52
$ws.Range("B52").Value = 8743
$ws.Range("C52").Value = 22
$ws.Range("D52").Value = 5922
$ws.Range("E52").Value = 2515
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = 306

# Row 56
$ws.Range("B56").Value = 7728
$ws.Range("C56").Value = 186
$ws.Range("D56").Value = 4062
$ws.Range("E56").Value = 3091
$ws.Range("G56").Value = 7
$ws.Range("H56").Value = 575

# Row 61 (now Moldavia)
$ws.Range("B61").Value = 6704
$ws.Range("C61").Value = 151
$ws.Range("D61").Value = 2953
$ws.Range("E61").Value = 3518
$ws.Range("G61").Value = 5
$ws.Range("H61").Value = 233

# Row 62 (now Nigeria, keeps old Nigeria figures)
$ws.Range("B62").Value = 6677
$ws.Range("D62").Value = 1840
$ws.Range("E62").Value = 4637
$ws.Range("H62").Value = 200

# Row 68
$ws.Range("B68").Value = 3980
$ws.Range("C68").Value = 9
$ws.Range("D68").Value = 3741
$ws.Range("E68").Value = 130

# Row 78
$ws.Range("B78").Value = 2853
$ws.Range("C78").Value = 3
$ws.Range("E78").Value = 1311
$ws.Range("G78").Value = 2
$ws.Range("H78").Value = 168

# Row 90
$ws.Range("D90").Value = 1790
$ws.Range("E90").Value = 3

# Row 111
$ws.Range("B111").Value = 923
$ws.Range("C111").Value = 1
$ws.Range("E111").Value = 390

# Row 132
$ws.Range("B132").Value = 457
$ws.Range("C132").Value = 30
$ws.Range("E132").Value = 405

# Row 142
$ws.Range("D142").Value = 303
$ws.Range("E142").Value = 9

# Row 151
$ws.Range("B151").Value = 220
$ws.Range("C151").Value = 3
$ws.Range("D151").Value = 112
$ws.Range("E151").Value = 106

# Row 163
$ws.Range("B163").Value = 135
$ws.Range("C163").Value = 5
$ws.Range("E163").Value = 75
$ws.Range("G163").Value = 1
$ws.Range("H163").Value = 3
